# Helper: write a numeric-looking string (e.g. "1.003", "0.3424") into a cell
# as literal TEXT, not a Number. A plain ".Value = '1.003'" assignment lets
# Excel's type inference parse it into a float (losing the original digits,
# e.g. trailing zeros) and keep it numeric. Instead, stash the string via a
# text formula (="1.003"), copy the cell, and paste-special just the value
# back onto itself (xlPasteValues = -4163). Excel pastes the formula's
# *string result* as a literal value, preserving exact text and leaving the
# cell's number format/style untouched (no "Text" format / quote-prefix
# style gets attached, unlike NumberFormat="@" or a leading apostrophe).
function Set-TextValue($ws, $addr, $text) {
    $cell = $ws.Range($addr)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (prices + 1h volume %) from the latest GitHub Actions
# scrape run. Plain string values (URLs, coin names, "xx.xx%" volume deltas)
# are safe to assign directly since Excel can't parse them as numbers;
# price values that look like plain numbers go through Set-TextValue so they
# stay text, matching the sheet's existing inlineStr storage.

$ws.Range("D2").Value = "22.390.68"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.571.99"
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws "D5" "1.003"
$ws.Range("E5").Value = "  +0.16%  "
Set-TextValue $ws "D6" "291.56"
$ws.Range("E6").Value = "  +0.53%  "
Set-TextValue $ws "D7" "0.3762"
$ws.Range("E7").Value = "  +2.34%  "
Set-TextValue $ws "D8" "50.14"
Set-TextValue $ws "D9" "0.3424"
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("E10").Value = "  +0.64%  "
Set-TextValue $ws "D11" "1.150"
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("E13").Value = "  -0.25%  "
Set-TextValue $ws "D14" "6.024"
$ws.Range("E14").Value = "  -0.45%  "
Set-TextValue $ws "D15" "6.965"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").Value = "1.578.74"
$ws.Range("E16").Value = "  +0.63%  "
Set-TextValue $ws "D17" "0.00001133"
$ws.Range("E17").Value = "  -0.16%  "
Set-TextValue $ws "D18" "90.09"
$ws.Range("E18").Value = "  +1.12%  "
Set-TextValue $ws "D19" "0.06727"
$ws.Range("E19").Value = "  -0.39%  "
Set-TextValue $ws "D21" "16.77"
$ws.Range("E21").Value = "  +1.50%  "
Set-TextValue $ws "D22" "6.229"
$ws.Range("E22").Value = "  -0.06%  "
Set-TextValue $ws "D23" "12.02"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "22.405.34"
$ws.Range("E24").Value = "  +0.07%  "
Set-TextValue $ws "D25" "2.395"
$ws.Range("E25").Value = "  +0.82%  "
Set-TextValue $ws "D26" "2.667"
$ws.Range("E26").Value = "  -11.17%  "
$ws.Range("E27").Value = "  +1.49%  "
Set-TextValue $ws "D28" "146.94"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("E29").Value = "  +1.44%  "
Set-TextValue $ws "D30" "126.37"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").Value = "1.750.55"
$ws.Range("E31").Value = "  +0.26%  "
Set-TextValue $ws "D32" "6.154"
$ws.Range("E32").Value = "  -1.88%  "
Set-TextValue $ws "D33" "1.997"
$ws.Range("E33").Value = "  -0.24%  "
Set-TextValue $ws "D34" "0.9837"
$ws.Range("E34").Value = "  -5.82%  "
Set-TextValue $ws "D35" "9.996"
$ws.Range("E35").Value = "  -3.28%  "
Set-TextValue $ws "D36" "0.08510"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws "D37" "1.396"
$ws.Range("E37").Value = "  +11.80%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D38" "0.02534"
$ws.Range("E38").Value = "  -0.41%  "
Set-TextValue $ws "D39" "0.2315"
$ws.Range("E39").Value = "  -0.92%  "
Set-TextValue $ws "D40" "0.06579"
$ws.Range("E40").Value = "  +0.76%  "
Set-TextValue $ws "D41" "5.416"
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws "D42" "0.6406"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws "D43" "11.46"
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("E44").Value = "  +0.23%  "
Set-TextValue $ws "D45" "14.05"
$ws.Range("E45").Value = "  -3.28%  "
Set-TextValue $ws "D46" "3.800"
$ws.Range("E46").Value = "  +0.81%  "
Set-TextValue $ws "D47" "0.5973"
$ws.Range("E47").Value = "  -0.41%  "
Set-TextValue $ws "D48" "1.293"
$ws.Range("E48").Value = "  +1.64%  "
Set-TextValue $ws "D49" "2.090"
$ws.Range("E49").Value = "  -1.86%  "
Set-TextValue $ws "D50" "125.48"
$ws.Range("E50").Value = "  +1.01%  "
Set-TextValue $ws "D51" "0.07331"
$ws.Range("E51").Value = "  +0.67%  "

# release the clipboard / marching-ants selection left by Copy()
$excel.CutCopyMode = $false
